# Apply the "Trade #3 closed" update to the live trading results workbook.
#
# Summary of changes:
#  - Summary sheet: update running totals now that a 3rd (losing) trade closed
#  - Strategy Status sheet: update the MarketMaking strategy row with new totals
#  - All Trades sheet: append the new trade row (#3)
#  - MarketMaking sheet: append the same new trade row (#3)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1500.09   # Current Capital
$summary.Range("B4").Value = 0.09      # Total P&L $
$summary.Range("B5").Value = 0.6       # Total P&L %
$summary.Range("B6").Value = 3         # Total Trades
$summary.Range("B8").Value = 2         # Losing Trades
$summary.Range("B9").Value = 33.33     # Win Rate %

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 100.09     # Capital
$status.Range("D6").Value = 3          # Trades
$status.Range("E6").Value = 0.09       # P&L $
$status.Range("F6").Value = 0.09       # P&L %
$status.Range("G6").Value = 33.33      # Win Rate %

# ---------------------------------------------------------------------------
# Helper: write the new Trade #3 row into a trades-style worksheet (used by
# both "All Trades" and "MarketMaking" sheets, which share the same layout).
# ---------------------------------------------------------------------------
function Add-Trade3Row($ws) {
    $ws.Range("A4").Value = 3
    # Force the date/time columns to stay plain text instead of being
    # auto-converted to Excel date/time serial values.
    $ws.Range("B4").Value = "'2026-02-17"
    $ws.Range("B4").Style = "Normal"
    $ws.Range("C4").Value = "23:52:01"
    $ws.Range("D4").Value = "MarketMaking"
    $ws.Range("E4").Value = "DOWN"
    $ws.Range("F4").Value = 0.1
    $ws.Range("G4").Value = 0.093361
    $ws.Range("H4").Value = "CLOSED"
    $ws.Range("I4").Value = -6.639
    $ws.Range("J4").Value = -0.01
    $ws.Range("K4").Value = 100.09
    $ws.Range("L4").Value = 0
    $ws.Range("M4").Value = 0
    $ws.Range("N4").Value = 0.6
    $ws.Range("O4").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P4").Value = "early_exit"
    $ws.Range("Q4").Value = 0.14
}

# ---------------------------------------------------------------------------
# 3. All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade3Row $allTrades

# ---------------------------------------------------------------------------
# 4. MarketMaking sheet
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade3Row $marketMaking
